$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 'Starting SoC (%)'
$ws.Cells.Item(6, 2).Value = 97
$ws.Cells.Item(7, 1).Value = 'Ending SoC (%)'
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(8, 1).Value = 'Total distance covered (km)'
$ws.Cells.Item(8, 2).Value = 69.78145244294669
$ws.Cells.Item(9, 1).Value = 'Total energy consumption(WH/KM)'
$ws.Cells.Item(9, 2).Value = 26.41004406635532
$ws.Cells.Item(10, 1).Value = 'Total SOC consumed(%)'
$ws.Cells.Item(10, 2).Value = 97
$ws.Cells.Item(12, 1).Value = 'Peak Power(kW)'
$ws.Cells.Item(12, 2).Value = 2438.510646
$ws.Cells.Item(13, 1).Value = 'Average Power(kW)'
$ws.Cells.Item(13, 2).Value = -813.6561739657836
$ws.Cells.Item(14, 1).Value = 'Total Energy Regenerated(kWh)'
$ws.Cells.Item(14, 2).Value = 22.40473596861111
$ws.Cells.Item(15, 1).Value = 'Regenerative Effectiveness(kWh)'
$ws.Cells.Item(15, 2).Value = -1.201109951715442
$ws.Cells.Item(16, 1).Value = 'Highest Cell Voltage(V)'
$ws.Cells.Item(16, 2).Value = 3.339
$ws.Cells.Item(17, 1).Value = 'Lowest Cell Voltage(V)'
$ws.Cells.Item(17, 2).Value = 2.985
$ws.Cells.Item(18, 1).Value = 'Difference in Cell Voltage(V)'
$ws.Cells.Item(18, 2).Value = 0.3540000000000001
$ws.Cells.Item(19, 1).Value = 'Minimum Temperature(C)'
$ws.Cells.Item(19, 2).Value = 23
$ws.Cells.Item(20, 1).Value = 'Maximum Temperature(C)'
$ws.Cells.Item(20, 2).Value = 44
$ws.Cells.Item(21, 1).Value = 'Difference in Temperature(C)'
$ws.Cells.Item(21, 2).Value = 21
$ws.Cells.Item(22, 1).Value = 'Maximum Fet Temperature-BMS(C)'
$ws.Cells.Item(22, 2).Value = 58
$ws.Cells.Item(23, 1).Value = 'Maximum Afe Temperature-BMS(C)'
$ws.Cells.Item(23, 2).Value = 62
$ws.Cells.Item(24, 1).Value = 'Maximum PCB Temperature-BMS(C)'
$ws.Cells.Item(24, 2).Value = 57
$ws.Cells.Item(25, 1).Value = 'Maximum MCU Temperature(C)'
$ws.Cells.Item(25, 2).Value = 42
$ws.Cells.Item(26, 1).Value = 'Maximum Motor Temperature(C)'
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(27, 1).Value = 'Abnormal Motor Temperature Detected(C)'
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(28, 1).Value = 'highest cell temp(C)'
$ws.Cells.Item(28, 2).Value = 44
$ws.Cells.Item(29, 1).Value = 'lowest cell temp(C)'
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(30, 1).Value = 'Difference between Highest and Lowest Cell Temperature at 100% SOC(C)'
$ws.Cells.Item(30, 2).Value = 44
$ws.Cells.Item(31, 1).Value = 'Battery Voltage(V)'
$ws.Cells.Item(31, 2).Value = 55
$ws.Cells.Item(32, 1).Value = 'Total energy charged(kWh)'
$ws.Cells.Item(32, 2).Value = 1.970930148611111
$ws.Cells.Item(33, 1).Value = 'Electricity consumption units(kW)'
$ws.Cells.Item(33, 2).Value = 0.00000006687194293836812
$ws.Cells.Item(34, 1).Value = 'Idling time percentage'
$ws.Cells.Item(34, 2).Value = 18.29305715381665
$ws.Cells.Item(35, 1).Value = 'Time spent in 0-10 km/h'
$ws.Cells.Item(35, 2).Value = 4.760260836210203
$ws.Cells.Item(36, 1).Value = 'Time spent in 10-20 km/h'
$ws.Cells.Item(36, 2).Value = 1.663470144482803
$ws.Cells.Item(37, 1).Value = 'Time spent in 20-30 km/h'
$ws.Cells.Item(37, 2).Value = 3.11724843370413
$ws.Cells.Item(38, 1).Value = 'Time spent in 30-40 km/h'
$ws.Cells.Item(38, 2).Value = 41.70950006393044
$ws.Cells.Item(39, 1).Value = 'Time spent in 40-50 km/h'
$ws.Cells.Item(39, 2).Value = 28.17286791970336
$ws.Cells.Item(40, 1).Value = 'Time spent in 50-60 km/h'
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(41, 1).Value = 'Time spent in 60-70 km/h'
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(42, 1).Value = 'Time spent in 70-80 km/h'
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(43, 1).Value = 'Time spent in 80-90 km/h'
$ws.Cells.Item(43, 2).Value = 0
